# Updates Leve profit-tracking values across the Sheets workbook (scheduled price refresh).
# Each block edits currentAveragePrice / LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ
# columns (H,I,J,K,L,M,N) for specific rows per crafting-class worksheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 186
$ws.Range("I33").Value = 186
$ws.Range("K33").Value = 186
$ws.Range("M33").Value = 43

# Row 43
$ws.Range("H43").Value = 3999
$ws.Range("J43").Value = 3999
$ws.Range("L43").Value = 3999
$ws.Range("N43").Value = -4137

# Row 62
$ws.Range("H62").Value = 2607.8333
$ws.Range("I62").Value = 1899
$ws.Range("J62").Value = 3316.6667
$ws.Range("K62").Value = 1899
$ws.Range("L62").Value = 3316.6667
$ws.Range("M62").Value = -1275
$ws.Range("N62").Value = -4564.6667

# Row 65
$ws.Range("H65").Value = 2607.8333
$ws.Range("I65").Value = 1899
$ws.Range("J65").Value = 3316.6667
$ws.Range("K65").Value = 9495
$ws.Range("L65").Value = 16583.3335
$ws.Range("M65").Value = -6375
$ws.Range("N65").Value = -22823.3335

# Row 103
$ws.Range("H103").Value = 859.625
$ws.Range("J103").Value = 859.625
$ws.Range("L103").Value = 2578.875
$ws.Range("N103").Value = -3750.875

# Row 138
$ws.Range("H138").Value = 2301.75
$ws.Range("I138").Value = 803.5
$ws.Range("J138").Value = 3800
$ws.Range("K138").Value = 2410.5
$ws.Range("L138").Value = 11400
$ws.Range("M138").Value = 2729.5
$ws.Range("N138").Value = -21680

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3918.4102
$ws.Range("J32").Value = 3300
$ws.Range("L32").Value = 3300
$ws.Range("N32").Value = -3874

# Row 61
$ws.Range("H61").Value = 3951.6667
$ws.Range("I61").Value = 3951.6667
$ws.Range("K61").Value = 3951.6667
$ws.Range("M61").Value = -3739.6667

# Row 74
$ws.Range("H74").Value = 5231.769
$ws.Range("I74").Value = 4973.909
$ws.Range("J74").Value = 6650
$ws.Range("K74").Value = 4973.909
$ws.Range("L74").Value = 6650
$ws.Range("M74").Value = -4099.909
$ws.Range("N74").Value = -8398

# Row 77
$ws.Range("H77").Value = 5231.769
$ws.Range("I77").Value = 4973.909
$ws.Range("J77").Value = 6650
$ws.Range("K77").Value = 24869.545
$ws.Range("L77").Value = 33250
$ws.Range("M77").Value = -20501.545
$ws.Range("N77").Value = -41986

# Row 136
$ws.Range("H136").Value = 3951.6667
$ws.Range("I136").Value = 3951.6667
$ws.Range("K136").Value = 11855.0001
$ws.Range("M136").Value = -9305.000100000001

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 744.6667
$ws.Range("I20").Value = 692.25
$ws.Range("K20").Value = 692.25
$ws.Range("M20").Value = -445.25

# Row 62
$ws.Range("H62").Value = 35000
$ws.Range("J62").Value = 35000
$ws.Range("L62").Value = 35000
$ws.Range("N62").Value = -36372

# Row 65
$ws.Range("H65").Value = 35000
$ws.Range("J65").Value = 35000
$ws.Range("L65").Value = 105000
$ws.Range("N65").Value = -111864

# Row 106
$ws.Range("H106").Value = 21511.4
$ws.Range("J106").Value = 21511.4
$ws.Range("L106").Value = 21511.4
$ws.Range("N106").Value = -24035.4

# Row 134
$ws.Range("H134").Value = 6716.5
$ws.Range("I134").Value = 6252.5835
$ws.Range("K134").Value = 18757.7505
$ws.Range("M134").Value = -16222.7505

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3449.4119
$ws.Range("I31").Value = 2589.6
$ws.Range("K31").Value = 2589.6
$ws.Range("M31").Value = -2294.6

# Row 34
$ws.Range("H34").Value = 3449.4119
$ws.Range("I34").Value = 2589.6
$ws.Range("K34").Value = 2589.6
$ws.Range("M34").Value = -2387.6

# Row 99
$ws.Range("H99").Value = 4078.5715
$ws.Range("I99").Value = 5325
$ws.Range("J99").Value = 2416.6667
$ws.Range("K99").Value = 5325
$ws.Range("L99").Value = 2416.6667
$ws.Range("M99").Value = -3827
$ws.Range("N99").Value = -5412.6667

# Row 126
$ws.Range("H126").Value = 4078.5715
$ws.Range("I126").Value = 5325
$ws.Range("J126").Value = 2416.6667
$ws.Range("K126").Value = 15975
$ws.Range("L126").Value = 7250.000100000001
$ws.Range("M126").Value = -13505
$ws.Range("N126").Value = -12190.0001

# Row 134
$ws.Range("H134").Value = 3199
$ws.Range("I134").Value = 3149.8462
$ws.Range("K134").Value = 9449.5386
$ws.Range("M134").Value = -6914.5386

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 32774.438
$ws.Range("I4").Value = 51558.9
$ws.Range("J4").Value = 1467
$ws.Range("K4").Value = 154676.7
$ws.Range("L4").Value = 4401
$ws.Range("M4").Value = -154564.7
$ws.Range("N4").Value = -4625

# Row 68
$ws.Range("H68").Value = 1998.25
$ws.Range("J68").Value = 1997.6666
$ws.Range("L68").Value = 5992.9998
$ws.Range("N68").Value = -7614.9998

# Row 71
$ws.Range("H71").Value = 1998.25
$ws.Range("J71").Value = 1997.6666
$ws.Range("L71").Value = 17978.9994
$ws.Range("N71").Value = -26090.9994

# Row 86
$ws.Range("H86").Value = 841.4
$ws.Range("I86").Value = 464.33334
$ws.Range("K86").Value = 1393.00002
$ws.Range("M86").Value = -207.0000199999999

# Row 89
$ws.Range("H89").Value = 841.4
$ws.Range("I89").Value = 464.33334
$ws.Range("K89").Value = 4179.00006
$ws.Range("M89").Value = 1748.99994

# Row 113
$ws.Range("H113").Value = 2255.25
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 2255.25
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 6765.75
$ws.Range("M113").Value = ""
$ws.Range("N113").Value = -11105.75

$ws = $wb.Worksheets.Item("GSM")
# Row 3
$ws.Range("H3").Value = 45455860
$ws.Range("I3").Value = 71428824
$ws.Range("K3").Value = 71428824
$ws.Range("M3").Value = -71428708

# Row 11
$ws.Range("H11").Value = 2637692.5
$ws.Range("I11").Value = 20003
$ws.Range("J11").Value = 2706579
$ws.Range("K11").Value = 20003
$ws.Range("L11").Value = 2706579
$ws.Range("M11").Value = -19864
$ws.Range("N11").Value = -2706857

# Row 44
$ws.Range("H44").Value = 30001
$ws.Range("J44").Value = 30001
$ws.Range("L44").Value = 30001
$ws.Range("N44").Value = -31193

# Row 122
$ws.Range("H122").Value = 13891912
$ws.Range("I122").Value = 17859028
$ws.Range("K122").Value = 53577084
$ws.Range("M122").Value = -53574634

# Row 132
$ws.Range("H132").Value = 1733.6
$ws.Range("I132").Value = 1815.1111
$ws.Range("K132").Value = 5445.3333
$ws.Range("M132").Value = -2915.3333

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 4222.222
$ws.Range("I46").Value = 4000
$ws.Range("K46").Value = 4000
$ws.Range("M46").Value = -3812

# Row 55
$ws.Range("H55").Value = 427.85715
$ws.Range("I55").Value = 398.33334
$ws.Range("J55").Value = 450
$ws.Range("K55").Value = 398.33334
$ws.Range("L55").Value = 450
$ws.Range("M55").Value = -225.33334
$ws.Range("N55").Value = -796

# Row 132
$ws.Range("H132").Value = 5937.4443
$ws.Range("I132").Value = 4156.1665
$ws.Range("J132").Value = 9500
$ws.Range("K132").Value = 12468.4995
$ws.Range("L132").Value = 28500
$ws.Range("M132").Value = -9938.499500000002
$ws.Range("N132").Value = -33560

# Row 133
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").Value = ""
